# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-10-29 09:15:49
#
# Re-applies the upstream "Recorded By" ordering fix, refreshed Class
# Statistics figures, and newly-recorded sessions for groups B2D/B2E/B2F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) "Recorded By" (column G) name-ordering fixes.
#    "System, X[, system]" -> "X[, System], System" style re-orderings.
# ---------------------------------------------------------------------
$gFixes = @{
    2   = "System, system, backup@backdoor.com"
    4   = "backup@backdoor.com, System"
    7   = "admin@admin.com, System"
    11  = "dnasr281@gmail.com, System"
    17  = "dnasr281@gmail.com, System"
    29  = "System, system, backup@backdoor.com"
    31  = "backup@backdoor.com, System"
    34  = "admin@admin.com, System"
    38  = "dnasr281@gmail.com, System"
    44  = "dnasr281@gmail.com, System"
    56  = "System, system, backup@backdoor.com"
    58  = "backup@backdoor.com, System"
    61  = "admin@admin.com, System"
    65  = "dnasr281@gmail.com, System"
    71  = "dnasr281@gmail.com, System"
    96  = "dnasr281@gmail.com, System"
    97  = "dnasr281@gmail.com, System"
    99  = "dnasr281@gmail.com, System"
    122 = "dnasr281@gmail.com, System"
    123 = "dnasr281@gmail.com, System"
    125 = "dnasr281@gmail.com, System"
    148 = "dnasr281@gmail.com, System"
    149 = "dnasr281@gmail.com, System"
    151 = "dnasr281@gmail.com, System"
}

foreach ($row in $gFixes.Keys) {
    $ws.Range("G$row").Value = $gFixes[$row]
}

# ---------------------------------------------------------------------
# 2) Class Statistics (K/L columns) refresh.
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 138
$ws.Range("L7").Value = 0
$ws.Range("L9").Value = "86.8%"
$ws.Range("L10").Value = "71.1%"

# ---------------------------------------------------------------------
# 3) Per-group statistics rows (B2D/B2E/B2F) refresh.
# ---------------------------------------------------------------------
$ws.Range("O18").Value = 22
$ws.Range("P18").Value = 0
$ws.Range("R18").Value = "84.6%"
$ws.Range("S18").Value = "75.9%"

$ws.Range("O19").Value = 22
$ws.Range("P19").Value = 0
$ws.Range("R19").Value = "84.6%"
$ws.Range("S19").Value = "73.6%"

$ws.Range("O20").Value = 22
$ws.Range("P20").Value = 0
$ws.Range("R20").Value = "84.6%"
$ws.Range("S20").Value = "81.0%"

# ---------------------------------------------------------------------
# 4) Newly-recorded "session 22" rows for B2D (104), B2E (130), B2F (156).
#    They were "Not Recorded" (special red-ish style) and are now
#    "Recorded" like every other row, so re-stripe them with the normal
#    "Recorded" row format (copied from row 2) before filling in values.
# ---------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A104:I104").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A130:I130").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A156:I156").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("G104").Value = "dnasr281@gmail.com"
$ws.Range("H104").Value = "35/56"
$ws.Range("I104").Value = "Recorded"

$ws.Range("G130").Value = "dnasr281@gmail.com"
$ws.Range("H130").Value = "33/55"
$ws.Range("I130").Value = "Recorded"

$ws.Range("G156").Value = "dnasr281@gmail.com"
$ws.Range("H156").Value = "40/57"
$ws.Range("I156").Value = "Recorded"

# ---------------------------------------------------------------------
# 5) Column I width 14 -> 10 (match column H's width).
# ---------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(8).ColumnWidth()
